$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 updates (text/inline-string columns) - keep stored as text, matching
# the original formatting (e.g. "132.81000000" rather than numeric 132.81)
$textCells = @("E53", "F53", "H53", "J53", "K53")
$textValues = @{
    "E53" = "132.81000000"
    "F53" = "722431.00564000"
    "H53" = "90913639.86164160"
    "J53" = "376108.36843000"
    "K53" = "47396735.47413880"
}

foreach ($addr in $textCells) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $textValues[$addr]
    $r.Style = "Normal"
}

# Row 53 updates (numeric columns)
$ws.Range("I53").Value = 203945
$ws.Range("M53").Value = 132.8099999999998
$ws.Range("N53").Value = 127.565
$ws.Range("O53").Value = 127.3357142857143
$ws.Range("P53").Value = 139.2460000000001
$ws.Range("Q53").Value = 185.9640000000001
$ws.Range("R53").Value = 132.81
$ws.Range("S53").Value = 130.449175210461
$ws.Range("T53").Value = 138.8662249579151
$ws.Range("U53").Value = 164.7665530599285
$ws.Range("V53").Value = -25.90032810201339
$ws.Range("W53").Value = -26.34337349548186
$ws.Range("X53").Value = 0.4430453934684735
